$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 3
$ws.Range("B3").Value = "相对于1添加在熊市中相对稳定的价值类型因子"

$ws.Range("A4").Value = 4
$ws.Range("B4").Value = "相对于2添加了在熊市中相对稳定的价值类因子"

$ws.Range("B5").Select()
